$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.986.20"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.844.53"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.57%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.013"
$ws.Range("E5").Value = "  +0.45%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "309.03"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("E7").Value = "  +1.87%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3676"
$ws.Range("E8").Value = "  +1.38%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07223"
$ws.Range("E9").Value = "  +0.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9300"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "19.83"
$ws.Range("E11").Value = "  +1.19%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07764"
$ws.Range("E12").Value = "  +0.79%  "
$ws.Range("D13").Value = "1.840.07"
$ws.Range("E13").Value = "  -1.39%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.388"
$ws.Range("E14").Value = "  +1.94%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.466"
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.76"
$ws.Range("E16").Value = "  +0.44%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.017"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008656"
$ws.Range("E18").Value = "  +0.74%  "
$ws.Range("D20").Value = "26.999.05"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.53"
$ws.Range("E21").Value = "  +1.10%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.053"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.64"
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.922"
$ws.Range("E24").Value = "  -0.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.76"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("E26").Value = "  -0.33%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.987"
$ws.Range("E27").Value = "  -3.22%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "114.51"
$ws.Range("E28").Value = "  +0.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.950"
$ws.Range("E29").Value = "  +0.66%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.311"
$ws.Range("E31").Value = "  +3.78%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.177"
$ws.Range("E32").Value = "  +0.19%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.507"
$ws.Range("E33").Value = "  +0.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7358"
$ws.Range("E34").Value = "  -1.51%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.648"
$ws.Range("E35").Value = "  -8.07%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.113"
$ws.Range("E36").Value = "  +2.77%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01965"
$ws.Range("E37").Value = "  +1.23%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05255"
$ws.Range("E38").Value = "  +1.71%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.975"
$ws.Range("E39").Value = "  -0.29%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5246"
$ws.Range("E40").Value = "  +1.51%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.017"
$ws.Range("E41").Value = "  +1.56%  "
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.268"
$ws.Range("E43").Value = "  +1.37%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.56"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.4730"
$ws.Range("E45").Value = "  +0.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.015"
$ws.Range("E46").Value = "  +0.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.51"
$ws.Range("E47").Value = "  +0.80%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.606"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "65.42"
$ws.Range("E49").Value = "  +1.29%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06058"
$ws.Range("E50").Value = "  +0.28%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8914"
$ws.Range("E51").Value = "  +3.24%  "
